# "Generate Report for Handback"
#
# The localization-status report gets refreshed once translations have been
# handed back from the vendor: the Status column flips from "Ready for
# handoff" to "Handed back: in sync with en-US", and each locale sheet grows
# two new columns of data (Latest Target File / Latest Handback File) plus a
# real "Latest Handback DateTime" (replacing the 0001-01-01 00:00:00
# placeholder) for every row that was handed back.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBackStatus = "Handed back: in sync with en-US"

# The Overview sheet mirrors the same "Ready for handoff" status text for
# each locale column, so it flips to the handback wording too.
$overview.Range("B2").Value = $handedBackStatus
$overview.Range("C2").Value = $handedBackStatus
$overview.Range("B3").Value = $handedBackStatus
$overview.Range("C3").Value = $handedBackStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------

# Status column (C) for both data rows now reflects the handback.
$zhcn.Range("C2").Value = $handedBackStatus
$zhcn.Range("C3").Value = $handedBackStatus

# New column values for the two "handed back" rows.
$zhcn.Range("F2").Value = "0bf270da-3282-4773-9f61-7e4d661c7e0b.md"
$zhcn.Range("G2").Value = "0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-03-25 10:42:22"

$zhcn.Range("F3").Value = "3033e96b-7e9f-4711-b640-cd2f653ae591.md"
$zhcn.Range("G3").Value = "3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-03-25 10:42:22"

# Rebuild the hyperlinks collection so the new Target-File / Handback-File
# links land next to their row's existing Source/Handoff links, in the same
# relative order the report lists the columns (A, D, F, G per row).
$zhcn.Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/0bf270da-3282-4773-9f61-7e4d661c7e0b.md", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3220de6c95ad032a525d6e2420cb587a71961860/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.zh-cn.xlf", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/0bf270da-3282-4773-9f61-7e4d661c7e0b.md", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3220de6c95ad032a525d6e2420cb587a71961860/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.zh-cn.xlf", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.zh-cn.xlf") | Out-Null

$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/3033e96b-7e9f-4711-b640-cd2f653ae591.md", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3220de6c95ad032a525d6e2420cb587a71961860/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.zh-cn.xlf", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/3033e96b-7e9f-4711-b640-cd2f653ae591.md", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3220de6c95ad032a525d6e2420cb587a71961860/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.zh-cn.xlf", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------

$dede.Range("C2").Value = $handedBackStatus
$dede.Range("C3").Value = $handedBackStatus

# New column values for the two "handed back" rows.
$dede.Range("F2").Value = "0bf270da-3282-4773-9f61-7e4d661c7e0b.md"
$dede.Range("G2").Value = "0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.de-de.xlf"
$dede.Range("H2").Value = "2016-03-25 10:42:38"

$dede.Range("F3").Value = "3033e96b-7e9f-4711-b640-cd2f653ae591.md"
$dede.Range("G3").Value = "3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.de-de.xlf"
$dede.Range("H3").Value = "2016-03-25 10:42:38"

# Rebuild the hyperlinks collection in the same column order as zh-cn above.
$dede.Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/0bf270da-3282-4773-9f61-7e4d661c7e0b.md", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b9bb7830b02061688120bbf840846ee0713bbb7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.de-de.xlf", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/0bf270da-3282-4773-9f61-7e4d661c7e0b.md", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6b9bb7830b02061688120bbf840846ee0713bbb7/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.de-de.xlf", "", "", "0bf270da-3282-4773-9f61-7e4d661c7e0b.31acd5f75c47c7e850e1bc12fa1d2ac7d086920f.de-de.xlf") | Out-Null

$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/3033e96b-7e9f-4711-b640-cd2f653ae591.md", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b9bb7830b02061688120bbf840846ee0713bbb7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.de-de.xlf", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ee9b9e32f5231022d3b68e7033110ee3ab3af/e2e/3033e96b-7e9f-4711-b640-cd2f653ae591.md", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6b9bb7830b02061688120bbf840846ee0713bbb7/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.de-de.xlf", "", "", "3033e96b-7e9f-4711-b640-cd2f653ae591.de8a5bb5ef0845f615153588f0896370b8a2f23e.de-de.xlf") | Out-Null
